$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = "'69.206.36"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.74%  '

# Row 3
$ws.Range('D3').Value = "'3.517.30"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.15%  '

# Row 4
$ws.Range('D4').Value = "'1.00"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.03%  '

# Row 5
$ws.Range('D5').Value = "'569.95"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.86%  '

# Row 6
$ws.Range('D6').Value = "'182.73"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.65%  '

# Row 7
$ws.Range('D7').Value = "'3.508.55"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.16%  '

# Row 8
$ws.Range('D8').Value = "'0.613"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.45%  '

# Row 9
$ws.Range('E9').Value = '  +0.05%  '

# Row 10
$ws.Range('E10').Value = '  +6.55%  '

# Row 11
$ws.Range('D11').Value = "'0.637"
$ws.Range('D11').Style = 'Normal'

# Row 12
$ws.Range('D12').Value = "'53.88"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.47%  '

# Row 13
$ws.Range('E13').Value = '  +1.14%  '

# Row 14
$ws.Range('D14').Value = "'9.46"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.14%  '

# Row 15
$ws.Range('D15').Value = "'4.078.09"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.31%  '

# Row 16
$ws.Range('D16').Value = "'19.25"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.81%  '

# Row 17
$ws.Range('D17').Value = "'3.510.82"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.49%  '

# Row 18
$ws.Range('D18').Value = "'69.072.66"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.82%  '

# Row 19
$ws.Range('D19').Value = "'12.52"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.40%  '

# Row 20
$ws.Range('D20').Value = "'0.120"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.00%  '

# Row 21
$ws.Range('D21').Value = "'537.67"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +14.28%  '

# Row 23
$ws.Range('D23').Value = "'20.55"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +6.88%  '

# Row 24
$ws.Range('E24').Value = '  -1.35%  '

# Row 25
$ws.Range('D25').Value = "'4.37"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.87%  '

# Row 26
$ws.Range('D26').Value = "'93.79"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +6.42%  '

# Row 27
$ws.Range('E27').Value = '  +0.74%  '

# Row 28
$ws.Range('D28').Value = "'2.91"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.19%  '

# Row 29
$ws.Range('D29').Value = "'9.14"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.03%  '

# Row 30
$ws.Range('D30').Value = "'31.63"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.09%  '

# Row 31
$ws.Range('D31').Value = "'7.29"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.10%  '

# Row 32
$ws.Range('D32').Value = "'12.72"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.69%  '

# Row 33
$ws.Range('D33').Value = "'64.08"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.22%  '

# Row 34
$ws.Range('E34').Value = '  -4.22%  '

# Row 35
$ws.Range('D35').Value = "'568.99"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.16%  '

# Row 36
$ws.Range('B36').Value = 'InjectiveProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D36').Value = "'38.23"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.59%  '

# Row 37
$ws.Range('B37').Value = 'Fetch.AI'
$ws.Range('C37').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D37').Value = "'3.08"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +8.64%  '

# Row 38
$ws.Range('E38').Value = '  +0.06%  '

# Row 39
$ws.Range('D39').Value = "'0.398"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.91%  '

# Row 40
$ws.Range('D40').Value = "'0.0₃0762"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.25%  '

# Row 41
$ws.Range('E41').Value = '  -4.27%  '

# Row 42
$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D42').Value = "'3.06"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.40%  '

# Row 43
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').Value = "'3.34"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.96%  '

# Row 44
$ws.Range('B44').Value = 'ApeXProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D44').Value = "'3.50"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +5.56%  '

# Row 45
$ws.Range('B45').Value = 'ThetaToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D45').Value = "'2.97"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.49%  '

# Row 46
$ws.Range('D46').Value = "'0.0443"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.22%  '

# Row 47
$ws.Range('D47').Value = "'3.193.41"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.59%  '

# Row 48
$ws.Range('E48').Value = '  -2.30%  '

# Row 49
$ws.Range('D49').Value = "'0.134"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.96%  '

# Row 50
$ws.Range('D50').Value = "'0.998"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.21%  '

# Row 51
$ws.Range('D51').Value = "'136.87"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.35%  '
